# Daily attendance processing - 2026-01-22 15:47:56
# Swap the order of names in the "Recorded By" column cells that currently
# read "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# (Students/System co-signed a session, and the pair is now displayed with
# the System entry first.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
$lastCol = $usedRange.Columns.Count + $usedRange.Column - 1

# Locate the "Recorded By" column from the header row instead of hard-coding
# the index, in case the layout ever shifts.
$recordedByCol = 7
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
